$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: Insert a new "Meta description" paragraph right after the
# title (Heading1) paragraph, before the "Gameplay" Heading2 paragraph.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Style = "Normal"

$metaRange = $metaPara.Range
$metaStart = $metaRange.Start
$metaRange.Text = "Meta description: Play 6 Fruits for free and enjoy a retro-style slot game with simple and functional gameplay, 5 paylines, and solid payouts. Read our review for more information."

# Make "Meta description" (without the colon) bold.
$boldRange = $d.Range($metaStart, $metaStart + 16)
$boldRange.Font.Bold = 1

# ---------------------------------------------------------------------
# Change 2: Near the end of the document, remove the duplicated bold
# title paragraph and replace the text of the following italic
# paragraph with the new image-generation prompt.
# ---------------------------------------------------------------------

# Find the duplicated bold-title paragraph (skip the very first
# paragraph, which legitimately holds the document title) and delete
# the whole paragraph, including its paragraph mark.
for ($i = $d.Paragraphs.Count; $i -ge 2; $i--) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text.TrimEnd([char]13, [char]7)
    if ($text -eq "Play 6 Fruits for Free - Retro-Style Slot Review") {
        $killRange = $d.Range($para.Range.Start, $para.Range.End)
        $killRange.Delete()
        break
    }
}

# Replace the text of the final (italic) paragraph with the new prompt,
# keeping its existing run/paragraph formatting. The range excludes the
# trailing paragraph-mark character so the text is overwritten in place
# rather than inserted before the existing run.
for ($i = $d.Paragraphs.Count; $i -ge 2; $i--) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text.TrimEnd([char]13, [char]7)
    if ($text -eq "Play 6 Fruits for free and enjoy a retro-style slot game with simple and functional gameplay, 5 paylines, and solid payouts. Read our review for more information.") {
        $textRange = $d.Range($para.Range.Start, $para.Range.End - 1)
        $textRange.Text = "Prompt: Create a cartoon-style feature image for 6 Fruits, featuring a happy Maya warrior with glasses. The feature image should be vibrant and colorful, featuring a cartoon-style Maya warrior wearing glasses and looking excitedly at the 6 Fruits slot game. The warrior should be holding a handful of fruits, including watermelons, plums, lemons, grapes, and cherries, indicating the symbols present in the game. The background should depict a retro-style casino, with neon lights and shiny slot machines. The image should have an overall fun and playful vibe."
        break
    }
}
